# Repgrid funcional, procesado del xlsx
# Applies accent/typo corrections to shared text labels and a data fix in B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "Ángela"
$ws.Range("M1").Value = "Acéntón"

# B2 held the rating "3" as text; keep it textual while correcting it to "7".
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "7"
$ws.Range("B2").Style = "Normal"

$ws.Range("O13").Value = "Enérgico"
$ws.Range("A14").Value = "Egocéntrico"
$ws.Range("O14").Value = "Empático"
$ws.Range("A17").Value = "Frío"
$ws.Range("O17").Value = "Cálido"
$ws.Range("A19").Value = "Antipático"
$ws.Range("O20").Value = "Gruñón"
